$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in this sheet hold scraped numeric-looking strings (e.g. "233.30",
# "1.00") as literal TEXT. Excel auto-converts such strings to numbers on
# assignment (stripping formatting like trailing zeros), so each cell is
# briefly switched to a text number format for the write, then restored to
# General (its original format) so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.256.39"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("E2").NumberFormat = "General"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.357.69"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.70%  "
$ws.Range("E3").NumberFormat = "General"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E4").NumberFormat = "General"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.30"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("E5").NumberFormat = "General"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.05%  "
$ws.Range("E6").NumberFormat = "General"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.57"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +12.95%  "
$ws.Range("E7").NumberFormat = "General"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("E8").NumberFormat = "General"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.480"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.87%  "
$ws.Range("E9").NumberFormat = "General"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0976"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("E10").NumberFormat = "General"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.12"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("E11").NumberFormat = "General"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("B12").NumberFormat = "General"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C12").NumberFormat = "General"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.711.33"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.83%  "
$ws.Range("E12").NumberFormat = "General"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "TRON"
$ws.Range("B13").NumberFormat = "General"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E13").NumberFormat = "General"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.14"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.98%  "
$ws.Range("E14").NumberFormat = "General"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("E15").NumberFormat = "General"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.865"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("E16").NumberFormat = "General"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.352.91"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("E17").NumberFormat = "General"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.289.81"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E18").NumberFormat = "General"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000101"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.89%  "
$ws.Range("E19").NumberFormat = "General"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.34"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("E20").NumberFormat = "General"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "74.32"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("E21").NumberFormat = "General"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "249.73"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E22").NumberFormat = "General"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E23").NumberFormat = "General"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.78"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.63%  "
$ws.Range("E24").NumberFormat = "General"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.44"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E25").NumberFormat = "General"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("B26").NumberFormat = "General"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.01"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("E26").NumberFormat = "General"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Toncoin"
$ws.Range("B27").NumberFormat = "General"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C27").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.20"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.31%  "
$ws.Range("E27").NumberFormat = "General"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.39"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.50%  "
$ws.Range("E28").NumberFormat = "General"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.69"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E29").NumberFormat = "General"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.37%  "
$ws.Range("E30").NumberFormat = "General"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.130"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("E31").NumberFormat = "General"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("E32").NumberFormat = "General"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.46%  "
$ws.Range("E33").NumberFormat = "General"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E34").NumberFormat = "General"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.05"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E35").NumberFormat = "General"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.55"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.10%  "
$ws.Range("E36").NumberFormat = "General"

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "RenderToken"
$ws.Range("B37").NumberFormat = "General"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C37").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.69"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E37").NumberFormat = "General"

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("B38").NumberFormat = "General"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C38").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.43"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.36%  "
$ws.Range("E38").NumberFormat = "General"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("E39").NumberFormat = "General"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E40").NumberFormat = "General"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.92"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("E41").NumberFormat = "General"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "18.58"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +8.59%  "
$ws.Range("E42").NumberFormat = "General"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +8.33%  "
$ws.Range("E43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.20"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E44").NumberFormat = "General"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("E45").NumberFormat = "General"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.02%  "
$ws.Range("E46").NumberFormat = "General"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E47").NumberFormat = "General"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.442.56"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E48").NumberFormat = "General"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.579.75"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.82%  "
$ws.Range("E49").NumberFormat = "General"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.76"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E50").NumberFormat = "General"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.96%  "
$ws.Range("E51").NumberFormat = "General"
